# The "Prix Spot" sheet has two blocks of daily columns glued together:
# a short "XX-nov" block (…, 09-nov, 10-nov) immediately followed by a
# full "01-oct." .. "31-oct." block. The commit adds the missing
# "11-nov" day: a brand-new column is inserted right before the current
# "01-oct." column (DN), pushing the whole Oct. block one column to the
# right (DN:ER -> DO:ES) and growing the sheet dimension to A1:ES25.
#
# The new column has no data yet, so every data row just gets "-"
# (the same placeholder used by the other not-yet-available "xx-nov"
# columns DF:DM), except the header row, which gets the new day label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new blank column before column DN (column 118), shifting
# DN:ER one column to the right, i.e. to DO:ES.
$ws.Range("DN1").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Range("DN1").Value = "11-nov"

# Data rows 2-25: no data yet for this new day, same placeholder as the
# sibling "xx-nov" columns.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 118).Value = "-"
}
